# Reproduces the "Add files via upload" commit for Sample_P4Division.xlsx:
# new A/C operands, new "answer" (E) entries for the first few rows, refreshed
# G (expected-answer) text, and F formulas that now point at their own row
# instead of the row two above (E1/G1 -> E3/G3, etc).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A3").Value = 6
$ws.Range("C3").Value = 91
$ws.Range("E3").Value = "15 r2"
$ws.Range("F3").Formula = '=IF(E3="","Please answer",IF(E3<>G3,"Wrong","Correct"))'
$ws.Range("G3").Value = "15 r 1"

$ws.Range("A4").Value = 4
$ws.Range("C4").Value = 39
$ws.Range("E4").Value = "9 r 3"
$ws.Range("F4").Formula = '=IF(E4="","Please answer",IF(E4<>G4,"Wrong","Correct"))'
$ws.Range("G4").Value = "9 r 3"

$ws.Range("A5").Value = 6
$ws.Range("C5").Value = 33
$ws.Range("E5").Value = "5 r 3"
$ws.Range("F5").Formula = '=IF(E5="","Please answer",IF(E5<>G5,"Wrong","Correct"))'
$ws.Range("G5").Value = "5 r 3"

$ws.Range("A6").Value = 9
$ws.Range("C6").Value = 79
$ws.Range("F6").Formula = '=IF(E6="","Please answer",IF(E6<>G6,"Wrong","Correct"))'
$ws.Range("G6").Value = "8 r 7"

$ws.Range("A7").Value = 8
$ws.Range("C7").Value = 28
$ws.Range("F7").Formula = '=IF(E7="","Please answer",IF(E7<>G7,"Wrong","Correct"))'
$ws.Range("G7").Value = "3 r 4"

$ws.Range("A8").Value = 5
$ws.Range("C8").Value = 33
$ws.Range("F8").Formula = '=IF(E8="","Please answer",IF(E8<>G8,"Wrong","Correct"))'
$ws.Range("G8").Value = "6 r 3"

$ws.Range("A9").Value = 7
$ws.Range("C9").Value = 15
$ws.Range("F9").Formula = '=IF(E9="","Please answer",IF(E9<>G9,"Wrong","Correct"))'
$ws.Range("G9").Value = "2 r 1"

$ws.Range("A10").Value = 9
$ws.Range("C10").Value = 46
$ws.Range("F10").Formula = '=IF(E10="","Please answer",IF(E10<>G10,"Wrong","Correct"))'
$ws.Range("G10").Value = "5 r 1"

$ws.Range("A11").Value = 9
$ws.Range("C11").Value = 34
$ws.Range("F11").Formula = '=IF(E11="","Please answer",IF(E11<>G11,"Wrong","Correct"))'
$ws.Range("G11").Value = "3 r 7"

$ws.Range("A12").Value = 6
$ws.Range("C12").Value = 21
$ws.Range("F12").Formula = '=IF(E12="","Please answer",IF(E12<>G12,"Wrong","Correct"))'
$ws.Range("G12").Value = "3 r 3"

$ws.Range("A13").Value = 9
$ws.Range("C13").Value = 17
$ws.Range("F13").Formula = '=IF(E13="","Please answer",IF(E13<>G13,"Wrong","Correct"))'
$ws.Range("G13").Value = "1 r 8"

$ws.Range("A14").Value = 5
$ws.Range("C14").Value = 76
$ws.Range("F14").Formula = '=IF(E14="","Please answer",IF(E14<>G14,"Wrong","Correct"))'
$ws.Range("G14").Value = "15 r 1"

$ws.Range("A15").Value = 4
$ws.Range("C15").Value = 28
$ws.Range("F15").Formula = '=IF(E15="","Please answer",IF(E15<>G15,"Wrong","Correct"))'
$ws.Range("G15").Value = "7 r 0"

$ws.Range("A16").Value = 9
$ws.Range("C16").Value = 73
$ws.Range("F16").Formula = '=IF(E16="","Please answer",IF(E16<>G16,"Wrong","Correct"))'
$ws.Range("G16").Value = "8 r 1"

$ws.Range("A17").Value = 5
$ws.Range("C17").Value = 86
$ws.Range("F17").Formula = '=IF(E17="","Please answer",IF(E17<>G17,"Wrong","Correct"))'
$ws.Range("G17").Value = "17 r 1"

$ws.Range("A18").Value = 8
$ws.Range("C18").Value = 95
$ws.Range("F18").Formula = '=IF(E18="","Please answer",IF(E18<>G18,"Wrong","Correct"))'
$ws.Range("G18").Value = "11 r 7"

$ws.Range("A19").Value = 5
$ws.Range("C19").Value = 33
$ws.Range("F19").Formula = '=IF(E19="","Please answer",IF(E19<>G19,"Wrong","Correct"))'
$ws.Range("G19").Value = "6 r 3"

$ws.Range("A20").Value = 8
$ws.Range("C20").Value = 14
$ws.Range("F20").Formula = '=IF(E20="","Please answer",IF(E20<>G20,"Wrong","Correct"))'
$ws.Range("G20").Value = "1 r 6"

$ws.Range("A21").Value = 4
$ws.Range("C21").Value = 28
$ws.Range("F21").Formula = '=IF(E21="","Please answer",IF(E21<>G21,"Wrong","Correct"))'
$ws.Range("G21").Value = "7 r 0"

$ws.Range("A22").Value = 8
$ws.Range("C22").Value = 36
$ws.Range("F22").Formula = '=IF(E22="","Please answer",IF(E22<>G22,"Wrong","Correct"))'
$ws.Range("G22").Value = "4 r 4"

# Match the author's final selection (cell E6).
$ws.Range("E6").Select() | Out-Null
